# 30.06.2025 - maç sonuçları
# Add the missing match score (3-3) for the Çirihtalar vs Kural Kesiciler
# match on row 15, and move the active selection to H17 (as in the source
# workbook after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the score for the match in row 15 (Skor1 / Skor2 columns).
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3

# Update the saved selection/active cell to H17.
$ws.Range("H17").Select()
